$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 208063511
$ws.Range("B1").Value = "Omer"
$ws.Range("C1").Value = "Avisror"
$ws.Range("D1").Value = 12345
$ws.Range("E1").Value = $true

$ws.Range("A2").Value = 987654321
$ws.Range("B2").Value = "jimmy"
$ws.Range("C2").Value = "james"
$ws.Range("D2").Value = 54321
$ws.Range("E2").Value = $false

$ws.Range("A3").Value = 123456789
$ws.Range("B3").Value = "Peter"
$ws.Range("C3").Value = "Parker"
$ws.Range("D3").Value = "dibs"
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = 313301129
$ws.Range("B4").Value = "Liat"
$ws.Range("C4").Value = "Mulian"
$ws.Range("D4").Value = 546500
$ws.Range("E4").Value = $true

$ws.Columns.Item(1).ColumnWidth = 9.140625

$ws.Range("F10").Select()
